# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns to match
# the latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.148.43'
$ws.Range("E2").Value = '  +0.71%  '

# Row 3
$ws.Range("D3").Value = '1.680.47'
$ws.Range("E3").Value = '  +0.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"  # force text: "215.29" would otherwise parse as a number
$ws.Range("D5").Value = '215.29'
$ws.Range("E5").Value = '  +0.12%  '

# Row 6
$ws.Range("E6").Value = '  +0.27%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("E8").Value = '  +2.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"  # force text: "21.36" would otherwise parse as a number
$ws.Range("D9").Value = '21.36'
$ws.Range("E9").Value = '  +5.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"  # force text: "0.0888" would otherwise parse as a number
$ws.Range("D11").Value = '0.0888'
$ws.Range("E11").Value = '  +0.19%  '

# Row 12
$ws.Range("D12").Value = '1.916.85'
$ws.Range("E12").Value = '  +0.36%  '

# Row 13
$ws.Range("D13").Value = '1.678.41'
$ws.Range("E13").Value = '  +0.52%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"  # force text: "4.15" would otherwise parse as a number
$ws.Range("D14").Value = '4.15'
$ws.Range("E14").Value = '  +1.57%  '

# Row 15
$ws.Range("E15").Value = '  +2.09%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"  # force text: "66.29" would otherwise parse as a number
$ws.Range("D16").Value = '66.29'
$ws.Range("E16").Value = '  +0.99%  '

# Row 17
$ws.Range("D17").Value = '27.146.66'
$ws.Range("E17").Value = '  +0.69%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"  # force text: "239.36" would otherwise parse as a number
$ws.Range("D18").Value = '239.36'
$ws.Range("E18").Value = '  +1.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"  # force text: "8.08" would otherwise parse as a number
$ws.Range("D19").Value = '8.08'
$ws.Range("E19").Value = '  -0.44%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0743'
$ws.Range("E20").Value = '  +1.27%  '

# Row 21
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("E22").Value = '  +1.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"  # force text: "9.46" would otherwise parse as a number
$ws.Range("D23").Value = '9.46'
$ws.Range("E23").Value = '  +3.00%  '

# Row 24
$ws.Range("E24").Value = '  -2.73%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"  # force text: "147.11" would otherwise parse as a number
$ws.Range("D25").Value = '147.11'
$ws.Range("E25").Value = '  +1.12%  '

# Row 26
$ws.Range("E26").Value = '  +0.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"  # force text: "16.37" would otherwise parse as a number
$ws.Range("D27").Value = '16.37'
$ws.Range("E27").Value = '  +2.24%  '

# Row 28
$ws.Range("E28").Value = '  +0.40%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("E30").Value = '  +0.43%  '

# Row 31
$ws.Range("E31").Value = '  +0.15%  '

# Row 32
$ws.Range("D32").Value = '1.567.27'
$ws.Range("E32").Value = '  +6.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"  # force text: "3.37" would otherwise parse as a number
$ws.Range("D33").Value = '3.37'
$ws.Range("E33").Value = '  +1.34%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"  # force text: "3.21" would otherwise parse as a number
$ws.Range("D34").Value = '3.21'
$ws.Range("E34").Value = '  +2.80%  '

# Row 35
$ws.Range("E35").Value = '  +0.21%  '

# Row 36
$ws.Range("E36").Value = '  +2.43%  '

# Row 37
$ws.Range("E37").Value = '  -0.74%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"  # force text: "0.932" would otherwise parse as a number
$ws.Range("D38").Value = '0.932'
$ws.Range("E38").Value = '  +4.10%  '

# Row 39
$ws.Range("E39").Value = '  +2.27%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"  # force text: "1.06" would otherwise parse as a number
$ws.Range("D40").Value = '1.06'
$ws.Range("E40").Value = '  +1.72%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"  # force text: "69.18" would otherwise parse as a number
$ws.Range("D41").Value = '69.18'
$ws.Range("E41").Value = '  +2.84%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"  # force text: "5.56" would otherwise parse as a number
$ws.Range("D43").Value = '5.56'
$ws.Range("E43").Value = '  -5.08%  '

# Row 44
$ws.Range("E44").Value = '  -2.47%  '

# Row 45
$ws.Range("D45").Value = '1.826.06'
$ws.Range("E45").Value = '  +0.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"  # force text: "0.787" would otherwise parse as a number
$ws.Range("D46").Value = '0.787'
$ws.Range("E46").Value = '  +1.43%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"  # force text: "90.66" would otherwise parse as a number
$ws.Range("D47").Value = '90.66'
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("E48").Value = '  +3.57%  '

# Row 49
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  +0.66%  '

# Row 50
$ws.Range("E50").Value = '  +2.27%  '

# Row 51
$ws.Range("E51").Value = '  +5.17%  '
